$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''42.731.81'
$ws.Range("E2").Value = '''  -1.68%  '
$ws.Range("D3").Value = '''2.306.17'
$ws.Range("E3").Value = '''  -0.35%  '
$ws.Range("E4").Value = '''  -0.03%  '
$ws.Range("D5").Value = '''302.52'
$ws.Range("E5").Value = '''  -2.13%  '
$ws.Range("D6").Value = '''99.31'
$ws.Range("E6").Value = '''  -6.02%  '
$ws.Range("D7").Value = '''0.502'
$ws.Range("E7").Value = '''  -4.61%  '
$ws.Range("E8").Value = '''  -0.02%  '
$ws.Range("E9").Value = '''  -3.41%  '
$ws.Range("D10").Value = '''34.74'
$ws.Range("E10").Value = '''  -3.32%  '
$ws.Range("E11").Value = '''  -2.42%  '
$ws.Range("E12").Value = '''  +0.26%  '
$ws.Range("D13").Value = '''6.72'
$ws.Range("E13").Value = '''  -3.73%  '
$ws.Range("D14").Value = '''2.662.25'
$ws.Range("E14").Value = '''  -0.24%  '
$ws.Range("D16").Value = '''2.290.84'
$ws.Range("E16").Value = '''  -0.87%  '
$ws.Range("D17").Value = '''0.798'
$ws.Range("E17").Value = '''  -0.54%  '
$ws.Range("D18").Value = '''42.610.45'
$ws.Range("E18").Value = '''  -1.77%  '
$ws.Range("D19").Value = '''0.0₃0904'
$ws.Range("E19").Value = '''  -2.01%  '
$ws.Range("D20").Value = '''11.60'
$ws.Range("E20").Value = '''  -2.94%  '
$ws.Range("D21").Value = '''6.05'
$ws.Range("E21").Value = '''  -2.42%  '
$ws.Range("D22").Value = '''67.82'
$ws.Range("E22").Value = '''  -0.21%  '
$ws.Range("D23").Value = '''235.67'
$ws.Range("E23").Value = '''  -2.13%  '
$ws.Range("E24").Value = '''  -3.56%  '
$ws.Range("D25").Value = '''2.50'
$ws.Range("E25").Value = '''  -4.03%  '
$ws.Range("E26").Value = '''  -0.19%  '
$ws.Range("D27").Value = '''24.75'
$ws.Range("E27").Value = '''  -0.38%  '
$ws.Range("D28").Value = '''2.17'
$ws.Range("E28").Value = '''  +2.86%  '
$ws.Range("D29").Value = '''34.32'
$ws.Range("E29").Value = '''  -5.56%  '
$ws.Range("D30").Value = '''164.88'
$ws.Range("E30").Value = '''  +1.63%  '
$ws.Range("E31").Value = '''  -5.08%  '
$ws.Range("E32").Value = '''  +0.06%  '
$ws.Range("E33").Value = '''  -4.63%  '
$ws.Range("E34").Value = '''  -4.81%  '
$ws.Range("D35").Value = '''4.45'
$ws.Range("E35").Value = '''  -3.03%  '
$ws.Range("D36").Value = '''16.73'
$ws.Range("E36").Value = '''  -8.89%  '
$ws.Range("B37").Value = '''Hedera'
$ws.Range("C37").Value = '''https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D37").Value = '''0.0695'
$ws.Range("E37").Value = '''  -5.35%  '
$ws.Range("B38").Value = '''LidoDAOToken'
$ws.Range("C38").Value = '''https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D38").Value = '''2.89'
$ws.Range("E38").Value = '''  -4.00%  '
$ws.Range("D39").Value = '''1.80'
$ws.Range("E39").Value = '''  -3.65%  '
$ws.Range("D40").Value = '''0.100'
$ws.Range("E40").Value = '''  -5.24%  '
$ws.Range("D41").Value = '''0.111'
$ws.Range("E41").Value = '''  -3.37%  '
$ws.Range("D42").Value = '''2.47'
$ws.Range("E42").Value = '''  -1.13%  '
$ws.Range("D43").Value = '''1.964.53'
$ws.Range("E43").Value = '''  -0.12%  '
$ws.Range("E44").Value = '''  -3.53%  '
$ws.Range("B45").Value = '''FraxShare'
$ws.Range("C45").Value = '''https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").Value = '''10.21'
$ws.Range("E45").Value = '''  -0.41%  '
$ws.Range("B46").Value = '''EnergySwap'
$ws.Range("C46").Value = '''https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = '''18.25'
$ws.Range("E46").Value = '''  -2.62%  '
$ws.Range("E47").Value = '''  -6.79%  '
$ws.Range("D48").Value = '''55.42'
$ws.Range("E48").Value = '''  -4.57%  '
$ws.Range("D49").Value = '''2.527.69'
$ws.Range("E49").Value = '''  -0.36%  '
$ws.Range("D50").Value = '''2.83'
$ws.Range("E50").Value = '''  -3.47%  '
$ws.Range("E51").Value = '''  -0.53%  '
